$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")
$ws.Rows.Item(12).Delete()
$ws.Range("A13").Select() | Out-Null
